$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry captures a single cell whose text value changed per the diff.
# NumberFormat is forced to "@" (Text) before assignment so Excel's COM
# layer does not auto-coerce numeric-looking strings (prices) or percent-like
# strings (volumes, e.g. "0.91%") into Number/Percentage typed cells -- the
# source data stores these as plain text (inlineStr) values.
$updates = @(
    @{ Cell = "D2"; Value = "246.46" },
    @{ Cell = "E2"; Value = "0.91%" },
    @{ Cell = "D3"; Value = "29.78" },
    @{ Cell = "E3"; Value = "9.51%" },
    @{ Cell = "D4"; Value = "5.165" },
    @{ Cell = "E4"; Value = "1.29%" },
    @{ Cell = "D5"; Value = "0.05713" },
    @{ Cell = "E5"; Value = "0.97%" },
    @{ Cell = "D6"; Value = "6.609" },
    @{ Cell = "E6"; Value = "2.10%" },
    @{ Cell = "B7"; Value = "GateToken" },
    @{ Cell = "C7"; Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt" },
    @{ Cell = "D7"; Value = "3.052" },
    @{ Cell = "E7"; Value = "1.65%" },
    @{ Cell = "B8"; Value = "MXToken" },
    @{ Cell = "C8"; Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx" },
    @{ Cell = "D8"; Value = "0.8585" },
    @{ Cell = "E8"; Value = "4.67%" },
    @{ Cell = "B9"; Value = "FTXToken" },
    @{ Cell = "C9"; Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt" },
    @{ Cell = "D9"; Value = "0.8670" },
    @{ Cell = "E9"; Value = "2.98%" },
    @{ Cell = "B10"; Value = "WazirX" },
    @{ Cell = "C10"; Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx" },
    @{ Cell = "D10"; Value = "0.1362" },
    @{ Cell = "E10"; Value = "2.50%" },
    @{ Cell = "B11"; Value = "MandalaExchangeToken" },
    @{ Cell = "C11"; Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx" },
    @{ Cell = "D11"; Value = "0.07068" },
    @{ Cell = "E11"; Value = "2.04%" },
    @{ Cell = "B12"; Value = "BitrueCoin" },
    @{ Cell = "C12"; Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr" },
    @{ Cell = "D12"; Value = "0.02863" },
    @{ Cell = "E12"; Value = "-4.15%" },
    @{ Cell = "B13"; Value = "BitMartToken" },
    @{ Cell = "C13"; Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx" },
    @{ Cell = "D13"; Value = "0.09372" },
    @{ Cell = "B14"; Value = "BitForexToken" },
    @{ Cell = "C14"; Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf" },
    @{ Cell = "D14"; Value = "0.001525" },
    @{ Cell = "E14"; Value = "0.92%" },
    @{ Cell = "B15"; Value = "CoinExToken" },
    @{ Cell = "C15"; Value = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet" },
    @{ Cell = "D15"; Value = "0.04147" },
    @{ Cell = "E15"; Value = "-1.62%" },
    @{ Cell = "D16"; Value = "0.006182" },
    @{ Cell = "E16"; Value = "0.47%" },
    @{ Cell = "E17"; Value = "3,764.51%" },
    @{ Cell = "E18"; Value = "-1.01%" },
    @{ Cell = "B19"; Value = "BTSEToken" },
    @{ Cell = "C19"; Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse" },
    @{ Cell = "D19"; Value = "2.257" },
    @{ Cell = "E19"; Value = "-2.14%" },
    @{ Cell = "B20"; Value = "One" },
    @{ Cell = "C20"; Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one" },
    @{ Cell = "D20"; Value = "0.01026" },
    @{ Cell = "E20"; Value = "1,623.61%" },
    @{ Cell = "E21"; Value = "1.04%" },
    @{ Cell = "D22"; Value = "0.03262" },
    @{ Cell = "E22"; Value = "3.30%" },
    @{ Cell = "D23"; Value = "0.1301" },
    @{ Cell = "E23"; Value = "3.70%" },
    @{ Cell = "D24"; Value = "2.908" },
    @{ Cell = "E24"; Value = "-18.13%" },
    @{ Cell = "D25"; Value = "0.1381" },
    @{ Cell = "E25"; Value = "0.47%" },
    @{ Cell = "D26"; Value = "0.005091" },
    @{ Cell = "E26"; Value = "14.21%" },
    @{ Cell = "D27"; Value = "0.001223" },
    @{ Cell = "E27"; Value = "-0.31%" },
    @{ Cell = "E28"; Value = "23.52%" },
    @{ Cell = "D40"; Value = "0.03764" },
    @{ Cell = "E40"; Value = "2.95%" },
    @{ Cell = "B41"; Value = "KickToken" },
    @{ Cell = "C41"; Value = "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick" },
    @{ Cell = "D41"; Value = "0.005774" },
    @{ Cell = "E41"; Value = "-4.59%" },
    @{ Cell = "B42"; Value = "BKEXToken" },
    @{ Cell = "C42"; Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk" },
    @{ Cell = "D42"; Value = "0.1070" },
    @{ Cell = "E42"; Value = "1.70%" },
    @{ Cell = "B43"; Value = "CEJI" },
    @{ Cell = "C43"; Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji" },
    @{ Cell = "D43"; Value = "0.002101" },
    @{ Cell = "E43"; Value = "-8.67%" },
    @{ Cell = "D44"; Value = "0.009777" },
    @{ Cell = "E44"; Value = "18.02%" },
    @{ Cell = "D45"; Value = "0.00005113" },
    @{ Cell = "E45"; Value = "-3.86%" },
    @{ Cell = "E46"; Value = "0.03%" },
    @{ Cell = "D47"; Value = "0.07514" },
    @{ Cell = "E47"; Value = "-41.40%" },
    @{ Cell = "D48"; Value = "0.002736" },
    @{ Cell = "E48"; Value = "-2.95%" },
    @{ Cell = "E49"; Value = "0.03%" },
    @{ Cell = "D50"; Value = "0.0002001" },
    @{ Cell = "E50"; Value = "0.03%" }
)

foreach ($u in $updates) {
    $cell = $ws.Range($u.Cell)
    $cell.NumberFormat = "@"
    $cell.Value = $u.Value
}

Write-Host "Applied $($updates.Count) cell updates"
